# eelgrass_acres.xlsx update
# Commit: "eelgrass_acreage obtained from annual monitoring reports by the
#          Piscataqua Region Estuary Partnership"
#
# The sheet is renamed Sheet1 -> Sheet2 (and gets the next sheetId, 2) while the
# data table is replaced wholesale: headers change from
#   Year | Eelgrass_Coverage_Acres | Source
# to
#   Year | GB_ZM_acres | GBE_ZM_acres
# and the row set grows from 23 data rows (1996-2019) to 26 data rows (1996-2023),
# now carrying two numeric series per year instead of one value + a source link.
# The six newest rows (2016-2023) get an explicit "0" (integer) number format.

$wb = $excel.ActiveWorkbook

# Duplicate the existing sheet onto itself: the copy inherits sheetFormatPr (e.g.
# defaultRowHeight 14.5) from the original and is assigned the next sheetId (2).
# Then drop the original "Sheet1" (sheetId 1), leaving only the copy, which we
# rename to "Sheet2". Net result in xl/workbook.xml:
#   <sheet name="Sheet2" sheetId="2" r:id="rId1"/>
$oldSheet = $wb.Worksheets.Item(1)
$oldSheet.Copy($oldSheet) | Out-Null
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null
$ws = $wb.ActiveSheet
$ws.Name = "Sheet2"

# Wipe the old table (headers + 23 rows of old eelgrass data) before laying down
# the new one.
$ws.UsedRange.Clear()

# --- Header row ---
$ws.Cells.Item(1,1).Value = "Year"
$ws.Cells.Item(1,2).Value = "GB_ZM_acres"
$ws.Cells.Item(1,3).Value = "GBE_ZM_acres"

# --- Data rows: Year, GB_ZM_acres, GBE_ZM_acres (most recent year first) ---
$ws.Cells.Item(2,1).Value = 2023
$ws.Cells.Item(2,2).Value = 855.56
$ws.Cells.Item(2,2).NumberFormat = "0"
$ws.Cells.Item(2,3).Value = 1024.51
$ws.Cells.Item(2,3).NumberFormat = "0"

$ws.Cells.Item(3,1).Value = 2022
$ws.Cells.Item(3,2).Value = 1393.12
$ws.Cells.Item(3,2).NumberFormat = "0"
$ws.Cells.Item(3,3).Value = 1606.74
$ws.Cells.Item(3,3).NumberFormat = "0"

$ws.Cells.Item(4,1).Value = 2021
$ws.Cells.Item(4,2).Value = 1266.03
$ws.Cells.Item(4,2).NumberFormat = "0"
$ws.Cells.Item(4,3).Value = 1566.47
$ws.Cells.Item(4,3).NumberFormat = "0"

$ws.Cells.Item(5,1).Value = 2019
$ws.Cells.Item(5,2).Value = 1344.99
$ws.Cells.Item(5,2).NumberFormat = "0"
$ws.Cells.Item(5,3).Value = 1570.87
$ws.Cells.Item(5,3).NumberFormat = "0"

$ws.Cells.Item(6,1).Value = 2017
$ws.Cells.Item(6,2).Value = 1362.42
$ws.Cells.Item(6,2).NumberFormat = "0"
$ws.Cells.Item(6,3).Value = 1546.66
$ws.Cells.Item(6,3).NumberFormat = "0"

$ws.Cells.Item(7,1).Value = 2016
$ws.Cells.Item(7,2).Value = 1489.9
$ws.Cells.Item(7,2).NumberFormat = "0"
$ws.Cells.Item(7,3).Value = 1688.71
$ws.Cells.Item(7,3).NumberFormat = "0"

$ws.Cells.Item(8,1).Value = 2015
$ws.Cells.Item(8,2).Value = 1319
$ws.Cells.Item(8,3).Value = 1493

$ws.Cells.Item(9,1).Value = 2014
$ws.Cells.Item(9,2).Value = 1466
$ws.Cells.Item(9,3).Value = 1620

$ws.Cells.Item(10,1).Value = 2013
$ws.Cells.Item(10,2).Value = 1266
$ws.Cells.Item(10,3).Value = 1448

$ws.Cells.Item(11,1).Value = 2012
$ws.Cells.Item(11,2).Value = 1599
$ws.Cells.Item(11,3).Value = 1813

$ws.Cells.Item(12,1).Value = 2011
$ws.Cells.Item(12,2).Value = 1624
$ws.Cells.Item(12,3).Value = 1836

$ws.Cells.Item(13,1).Value = 2010
$ws.Cells.Item(13,2).Value = 1722
$ws.Cells.Item(13,3).Value = 1895

$ws.Cells.Item(14,1).Value = 2009
$ws.Cells.Item(14,2).Value = 1701
$ws.Cells.Item(14,3).Value = 1890

$ws.Cells.Item(15,1).Value = 2008
$ws.Cells.Item(15,2).Value = 1395
$ws.Cells.Item(15,3).Value = 1619

$ws.Cells.Item(16,1).Value = 2007
$ws.Cells.Item(16,2).Value = 1245
$ws.Cells.Item(16,3).Value = 1489

$ws.Cells.Item(17,1).Value = 2006
$ws.Cells.Item(17,2).Value = 1321
$ws.Cells.Item(17,3).Value = 1623

$ws.Cells.Item(18,1).Value = 2005
$ws.Cells.Item(18,2).Value = 2175
$ws.Cells.Item(18,3).Value = 2507

$ws.Cells.Item(19,1).Value = 2004
$ws.Cells.Item(19,2).Value = 2042
$ws.Cells.Item(19,3).Value = 2349

$ws.Cells.Item(20,1).Value = 2003
$ws.Cells.Item(20,2).Value = 1627
$ws.Cells.Item(20,3).Value = 1996

$ws.Cells.Item(21,1).Value = 2002
$ws.Cells.Item(21,2).Value = 1795
$ws.Cells.Item(21,3).Value = 2415

$ws.Cells.Item(22,1).Value = 2001
$ws.Cells.Item(22,2).Value = 2392
$ws.Cells.Item(22,3).Value = 2735

$ws.Cells.Item(23,1).Value = 2000
$ws.Cells.Item(23,2).Value = 1945
$ws.Cells.Item(23,3).Value = 2285

$ws.Cells.Item(24,1).Value = 1999
$ws.Cells.Item(24,2).Value = 2130
$ws.Cells.Item(24,3).Value = 2459

$ws.Cells.Item(25,1).Value = 1998
$ws.Cells.Item(25,2).Value = 2398

$ws.Cells.Item(26,1).Value = 1997
$ws.Cells.Item(26,2).Value = 2305

$ws.Cells.Item(27,1).Value = 1996
$ws.Cells.Item(27,2).Value = 2503
$ws.Cells.Item(27,3).Value = 2894

# --- Selection / view state left by the editor ---
$ws.Range("I34").Select() | Out-Null
